$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A1").Formula = "=MOD(12,5)"
$ws2.Range("B1").Value = "Jan"
$ws2.Range("C1").Value = "Feb"
$ws2.Range("D1").Value = "Mar"
$ws2.Range("F1").Value = "Jan"

$ws2.Columns.Item(8).Font.Color = 0
$ws2.Columns.Item(8).Delete()
